# Auto-generated BRVM data refresh (GitHub Actions MAJ)
$wb = $excel.ActiveWorkbook

# --- Sheet "Recommandations": update rows 2-51 (cols A-G) ---
$ws1 = $wb.Worksheets.Item("Recommandations")

# Row 2: BRVM - SERVICES PUBLICS
$ws1.Cells.Item(2, 1).Value = "BRVM - SERVICES PUBLICS"
$ws1.Cells.Item(2, 2).Value = 0
$ws1.Cells.Item(2, 3).Value = 10
$ws1.Cells.Item(2, 4).Value = 4292.9
$ws1.Cells.Item(2, 5).Value = 113.05
$ws1.Cells.Item(2, 6).Value = "🟡 Observer"
$ws1.Cells.Item(2, 7).Value = "➖ Neutre"

# Row 3: SAFCA CI
$ws1.Cells.Item(3, 1).Value = "SAFCA CI"
$ws1.Cells.Item(3, 2).Value = 0
$ws1.Cells.Item(3, 3).Value = 5
$ws1.Cells.Item(3, 4).Value = 3450
$ws1.Cells.Item(3, 5).Value = 690
$ws1.Cells.Item(3, 6).Value = "🟡 Observer"
$ws1.Cells.Item(3, 7).Value = "➖ Neutre"

# Row 4: CFAO MOTORS CI
$ws1.Cells.Item(4, 1).Value = "CFAO MOTORS CI"
$ws1.Cells.Item(4, 2).Value = 0
$ws1.Cells.Item(4, 3).Value = 5
$ws1.Cells.Item(4, 4).Value = 3375
$ws1.Cells.Item(4, 5).Value = 675
$ws1.Cells.Item(4, 6).Value = "🟡 Observer"
$ws1.Cells.Item(4, 7).Value = "➖ Neutre"

# Row 5: BRVM - AUTRES SECTEURS
$ws1.Cells.Item(5, 1).Value = "BRVM - AUTRES SECTEURS"
$ws1.Cells.Item(5, 2).Value = 0
$ws1.Cells.Item(5, 3).Value = 5
$ws1.Cells.Item(5, 4).Value = 3319.9
$ws1.Cells.Item(5, 5).Value = 665.88
$ws1.Cells.Item(5, 6).Value = "🟡 Observer"
$ws1.Cells.Item(5, 7).Value = "➖ Neutre"

# Row 6: UNIWAX CI
$ws1.Cells.Item(6, 1).Value = "UNIWAX CI"
$ws1.Cells.Item(6, 2).Value = 0
$ws1.Cells.Item(6, 3).Value = 5
$ws1.Cells.Item(6, 4).Value = 2970
$ws1.Cells.Item(6, 5).Value = 600
$ws1.Cells.Item(6, 6).Value = "🟡 Observer"
$ws1.Cells.Item(6, 7).Value = "➖ Neutre"

# Row 7: NEI-CEDA CI
$ws1.Cells.Item(7, 1).Value = "NEI-CEDA CI"
$ws1.Cells.Item(7, 2).Value = 0
$ws1.Cells.Item(7, 3).Value = 5
$ws1.Cells.Item(7, 4).Value = 2955
$ws1.Cells.Item(7, 5).Value = 595
$ws1.Cells.Item(7, 6).Value = "🟡 Observer"
$ws1.Cells.Item(7, 7).Value = "➖ Neutre"

# Row 8: SUCRIVOIRE
$ws1.Cells.Item(8, 1).Value = "SUCRIVOIRE"
$ws1.Cells.Item(8, 2).Value = 0
$ws1.Cells.Item(8, 3).Value = 3
$ws1.Cells.Item(8, 4).Value = 2945
$ws1.Cells.Item(8, 5).Value = 985
$ws1.Cells.Item(8, 6).Value = "🟡 Observer"
$ws1.Cells.Item(8, 7).Value = "➖ Neutre"

# Row 9: SETAO CI
$ws1.Cells.Item(9, 1).Value = "SETAO CI"
$ws1.Cells.Item(9, 2).Value = 0
$ws1.Cells.Item(9, 3).Value = 5
$ws1.Cells.Item(9, 4).Value = 2720
$ws1.Cells.Item(9, 5).Value = 575
$ws1.Cells.Item(9, 6).Value = "🟡 Observer"
$ws1.Cells.Item(9, 7).Value = "➖ Neutre"

# Row 10: AIR LIQUIDE CI
$ws1.Cells.Item(10, 1).Value = "AIR LIQUIDE CI"
$ws1.Cells.Item(10, 2).Value = 0
$ws1.Cells.Item(10, 3).Value = 5
$ws1.Cells.Item(10, 4).Value = 2710
$ws1.Cells.Item(10, 5).Value = 540
$ws1.Cells.Item(10, 6).Value = "🟡 Observer"
$ws1.Cells.Item(10, 7).Value = "➖ Neutre"

# Row 11: BRVM - TRANSPORT
$ws1.Cells.Item(11, 1).Value = "BRVM - TRANSPORT"
$ws1.Cells.Item(11, 2).Value = 0
$ws1.Cells.Item(11, 3).Value = 5
$ws1.Cells.Item(11, 4).Value = 1858.8
$ws1.Cells.Item(11, 5).Value = 348.8
$ws1.Cells.Item(11, 6).Value = "🟡 Observer"
$ws1.Cells.Item(11, 7).Value = "➖ Neutre"

# Row 12: BRVM - DISTRIBUTION
$ws1.Cells.Item(12, 1).Value = "BRVM - DISTRIBUTION"
$ws1.Cells.Item(12, 2).Value = 0
$ws1.Cells.Item(12, 3).Value = 5
$ws1.Cells.Item(12, 4).Value = 1841.82
$ws1.Cells.Item(12, 5).Value = 370.28
$ws1.Cells.Item(12, 6).Value = "🟡 Observer"
$ws1.Cells.Item(12, 7).Value = "➖ Neutre"

# Row 13: BRVM - AGRICULTURE
$ws1.Cells.Item(13, 1).Value = "BRVM - AGRICULTURE"
$ws1.Cells.Item(13, 2).Value = 0
$ws1.Cells.Item(13, 3).Value = 5
$ws1.Cells.Item(13, 4).Value = 1663.65
$ws1.Cells.Item(13, 5).Value = 333.59
$ws1.Cells.Item(13, 6).Value = "🟡 Observer"
$ws1.Cells.Item(13, 7).Value = "➖ Neutre"

# Row 14: BRVM - INDUSTRIE
$ws1.Cells.Item(14, 1).Value = "BRVM - INDUSTRIE"
$ws1.Cells.Item(14, 2).Value = 0
$ws1.Cells.Item(14, 3).Value = 5
$ws1.Cells.Item(14, 4).Value = 965.04
$ws1.Cells.Item(14, 5).Value = 192.93
$ws1.Cells.Item(14, 6).Value = "🟡 Observer"
$ws1.Cells.Item(14, 7).Value = "➖ Neutre"

# Row 15: BRVM-PRINCIPAL
$ws1.Cells.Item(15, 1).Value = "BRVM-PRINCIPAL"
$ws1.Cells.Item(15, 2).Value = 0
$ws1.Cells.Item(15, 3).Value = 5
$ws1.Cells.Item(15, 4).Value = 888.6799999999999
$ws1.Cells.Item(15, 5).Value = 176.76
$ws1.Cells.Item(15, 6).Value = "🟡 Observer"
$ws1.Cells.Item(15, 7).Value = "➖ Neutre"

# Row 16: BRVM - CONSOMMATION DE BASE
$ws1.Cells.Item(16, 1).Value = "BRVM - CONSOMMATION DE BASE"
$ws1.Cells.Item(16, 2).Value = 0
$ws1.Cells.Item(16, 3).Value = 5
$ws1.Cells.Item(16, 4).Value = 851.88
$ws1.Cells.Item(16, 5).Value = 171.42
$ws1.Cells.Item(16, 6).Value = "🟡 Observer"
$ws1.Cells.Item(16, 7).Value = "➖ Neutre"

# Row 17: BRVM - INDUSTRIELS
$ws1.Cells.Item(17, 1).Value = "BRVM - INDUSTRIELS"
$ws1.Cells.Item(17, 2).Value = 0
$ws1.Cells.Item(17, 3).Value = 5
$ws1.Cells.Item(17, 4).Value = 698.51
$ws1.Cells.Item(17, 5).Value = 129.78
$ws1.Cells.Item(17, 6).Value = "🟡 Observer"
$ws1.Cells.Item(17, 7).Value = "➖ Neutre"

# Row 18: BRVM-PRESTIGE
$ws1.Cells.Item(18, 1).Value = "BRVM-PRESTIGE"
$ws1.Cells.Item(18, 2).Value = 0
$ws1.Cells.Item(18, 3).Value = 5
$ws1.Cells.Item(18, 4).Value = 657.11
$ws1.Cells.Item(18, 5).Value = 131.19
$ws1.Cells.Item(18, 6).Value = "🟡 Observer"
$ws1.Cells.Item(18, 7).Value = "➖ Neutre"

# Row 19: BRVM - FINANCES
$ws1.Cells.Item(19, 1).Value = "BRVM - FINANCES"
$ws1.Cells.Item(19, 2).Value = 0
$ws1.Cells.Item(19, 3).Value = 5
$ws1.Cells.Item(19, 4).Value = 618.39
$ws1.Cells.Item(19, 5).Value = 123.18
$ws1.Cells.Item(19, 6).Value = "🟡 Observer"
$ws1.Cells.Item(19, 7).Value = "➖ Neutre"

# Row 20: BRVM - SERVICES FINANCIERS
$ws1.Cells.Item(20, 1).Value = "BRVM - SERVICES FINANCIERS"
$ws1.Cells.Item(20, 2).Value = 0
$ws1.Cells.Item(20, 3).Value = 5
$ws1.Cells.Item(20, 4).Value = 607.74
$ws1.Cells.Item(20, 5).Value = 121.06
$ws1.Cells.Item(20, 6).Value = "🟡 Observer"
$ws1.Cells.Item(20, 7).Value = "➖ Neutre"

# Row 21: BRVM - ENERGIE
$ws1.Cells.Item(21, 1).Value = "BRVM - ENERGIE"
$ws1.Cells.Item(21, 2).Value = 0
$ws1.Cells.Item(21, 3).Value = 5
$ws1.Cells.Item(21, 4).Value = 549.33
$ws1.Cells.Item(21, 5).Value = 110.27
$ws1.Cells.Item(21, 6).Value = "🟡 Observer"
$ws1.Cells.Item(21, 7).Value = "➖ Neutre"

# Row 22: BRVM - CONSOMMATION DISCRETIONNAIRE
$ws1.Cells.Item(22, 1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Cells.Item(22, 2).Value = 0
$ws1.Cells.Item(22, 3).Value = 5
$ws1.Cells.Item(22, 4).Value = 533.86
$ws1.Cells.Item(22, 5).Value = 107.22
$ws1.Cells.Item(22, 6).Value = "🟡 Observer"
$ws1.Cells.Item(22, 7).Value = "➖ Neutre"

# Row 23: BRVM - TELECOMMUNICATIONS
$ws1.Cells.Item(23, 1).Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Cells.Item(23, 2).Value = 0
$ws1.Cells.Item(23, 3).Value = 5
$ws1.Cells.Item(23, 4).Value = 484.9
$ws1.Cells.Item(23, 5).Value = 96.01000000000001
$ws1.Cells.Item(23, 6).Value = "🟡 Observer"
$ws1.Cells.Item(23, 7).Value = "➖ Neutre"

# Row 24: UNILEVER CI (UNLC)
$ws1.Cells.Item(24, 1).Value = "UNILEVER CI (UNLC)"
$ws1.Cells.Item(24, 2).Value = 4
$ws1.Cells.Item(24, 3).Value = 0
$ws1.Cells.Item(24, 4).Value = 29.92
$ws1.Cells.Item(24, 5).Value = 7.49
$ws1.Cells.Item(24, 6).Value = "🟢 Achat"
$ws1.Cells.Item(24, 7).Value = "✅ Renforcer"

# Row 25: NSIA BANQUE COTE D'IVOIRE (NSBC)
$ws1.Cells.Item(25, 1).Value = "NSIA BANQUE COTE D'IVOIRE (NSBC)"
$ws1.Cells.Item(25, 2).Value = 2
$ws1.Cells.Item(25, 3).Value = 0
$ws1.Cells.Item(25, 4).Value = 7.66
$ws1.Cells.Item(25, 5).Value = 3.03
$ws1.Cells.Item(25, 6).Value = "🟡 Observer"
$ws1.Cells.Item(25, 7).Value = "➖ Neutre"

# Row 26: SETAO CI (STAC)
$ws1.Cells.Item(26, 1).Value = "SETAO CI (STAC)"
$ws1.Cells.Item(26, 2).Value = 2
$ws1.Cells.Item(26, 3).Value = 2
$ws1.Cells.Item(26, 4).Value = 4.97
$ws1.Cells.Item(26, 5).Value = -2.59
$ws1.Cells.Item(26, 6).Value = "🟡 Observer"
$ws1.Cells.Item(26, 7).Value = "👀 À surveiller"

# Row 27: SAFCA CI (SAFC)
$ws1.Cells.Item(27, 1).Value = "SAFCA CI (SAFC)"
$ws1.Cells.Item(27, 2).Value = 1
$ws1.Cells.Item(27, 3).Value = 0
$ws1.Cells.Item(27, 4).Value = 4.55
$ws1.Cells.Item(27, 5).Value = 4.55
$ws1.Cells.Item(27, 6).Value = "🟡 Observer"
$ws1.Cells.Item(27, 7).Value = "➖ Neutre"

# Row 28: ECOBANK COTE D''IVOIRE (ECOC)
$ws1.Cells.Item(28, 1).Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws1.Cells.Item(28, 2).Value = 1
$ws1.Cells.Item(28, 3).Value = 0
$ws1.Cells.Item(28, 4).Value = 4.3
$ws1.Cells.Item(28, 5).Value = 4.3
$ws1.Cells.Item(28, 6).Value = "🟡 Observer"
$ws1.Cells.Item(28, 7).Value = "➖ Neutre"

# Row 29: CIE CI (CIEC)
$ws1.Cells.Item(29, 1).Value = "CIE CI (CIEC)"
$ws1.Cells.Item(29, 2).Value = 1
$ws1.Cells.Item(29, 3).Value = 0
$ws1.Cells.Item(29, 4).Value = 4
$ws1.Cells.Item(29, 5).Value = 4
$ws1.Cells.Item(29, 6).Value = "🟡 Observer"
$ws1.Cells.Item(29, 7).Value = "➖ Neutre"

# Row 30: TOTALENERGIES MARKETING CI (TTLC)
$ws1.Cells.Item(30, 1).Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Cells.Item(30, 2).Value = 1
$ws1.Cells.Item(30, 3).Value = 0
$ws1.Cells.Item(30, 4).Value = 3.39
$ws1.Cells.Item(30, 5).Value = 3.39
$ws1.Cells.Item(30, 6).Value = "🟡 Observer"
$ws1.Cells.Item(30, 7).Value = "➖ Neutre"

# Row 31: SUCRIVOIRE (SCRC)
$ws1.Cells.Item(31, 1).Value = "SUCRIVOIRE (SCRC)"
$ws1.Cells.Item(31, 2).Value = 1
$ws1.Cells.Item(31, 3).Value = 0
$ws1.Cells.Item(31, 4).Value = 3.09
$ws1.Cells.Item(31, 5).Value = 3.09
$ws1.Cells.Item(31, 6).Value = "🟡 Observer"
$ws1.Cells.Item(31, 7).Value = "➖ Neutre"

# Row 32: CFAO MOTORS CI (CFAC)
$ws1.Cells.Item(32, 1).Value = "CFAO MOTORS CI (CFAC)"
$ws1.Cells.Item(32, 2).Value = 1
$ws1.Cells.Item(32, 3).Value = 0
$ws1.Cells.Item(32, 4).Value = 3.03
$ws1.Cells.Item(32, 5).Value = 3.03
$ws1.Cells.Item(32, 6).Value = "🟡 Observer"
$ws1.Cells.Item(32, 7).Value = "➖ Neutre"

# Row 33: BERNABE CI (BNBC)
$ws1.Cells.Item(33, 1).Value = "BERNABE CI (BNBC)"
$ws1.Cells.Item(33, 2).Value = 1
$ws1.Cells.Item(33, 3).Value = 0
$ws1.Cells.Item(33, 4).Value = 3.02
$ws1.Cells.Item(33, 5).Value = 3.02
$ws1.Cells.Item(33, 6).Value = "🟡 Observer"
$ws1.Cells.Item(33, 7).Value = "➖ Neutre"

# Row 34: BANK OF AFRICA BN (BOAB)
$ws1.Cells.Item(34, 1).Value = "BANK OF AFRICA BN (BOAB)"
$ws1.Cells.Item(34, 2).Value = 1
$ws1.Cells.Item(34, 3).Value = 0
$ws1.Cells.Item(34, 4).Value = 2.81
$ws1.Cells.Item(34, 5).Value = 2.81
$ws1.Cells.Item(34, 6).Value = "🟡 Observer"
$ws1.Cells.Item(34, 7).Value = "➖ Neutre"

# Row 35: TOTAL
$ws1.Cells.Item(35, 1).Value = "TOTAL"
$ws1.Cells.Item(35, 2).Value = 0
$ws1.Cells.Item(35, 3).Value = 5
$ws1.Cells.Item(35, 4).Value = 0
$ws1.Cells.Item(35, 5).Value = 0
$ws1.Cells.Item(35, 6).Value = "🟡 Observer"
$ws1.Cells.Item(35, 7).Value = "➖ Neutre"

# Row 36: ONATEL BF (ONTBF)
$ws1.Cells.Item(36, 1).Value = "ONATEL BF (ONTBF)"
$ws1.Cells.Item(36, 2).Value = 1
$ws1.Cells.Item(36, 3).Value = 1
$ws1.Cells.Item(36, 4).Value = -0.02
$ws1.Cells.Item(36, 5).Value = 7.48
$ws1.Cells.Item(36, 6).Value = "🟡 Observer"
$ws1.Cells.Item(36, 7).Value = "👀 À surveiller"

# Row 37: TRACTAFRIC MOTORS CI (PRSC)
$ws1.Cells.Item(37, 1).Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Cells.Item(37, 2).Value = 1
$ws1.Cells.Item(37, 3).Value = 1
$ws1.Cells.Item(37, 4).Value = -0.02
$ws1.Cells.Item(37, 5).Value = 4.16
$ws1.Cells.Item(37, 6).Value = "🟡 Observer"
$ws1.Cells.Item(37, 7).Value = "👀 À surveiller"

# Row 38: ORANGE COTE D'IVOIRE (ORAC)
$ws1.Cells.Item(38, 1).Value = "ORANGE COTE D'IVOIRE (ORAC)"
$ws1.Cells.Item(38, 2).Value = 1
$ws1.Cells.Item(38, 3).Value = 1
$ws1.Cells.Item(38, 4).Value = -0.5600000000000001
$ws1.Cells.Item(38, 5).Value = 2.6
$ws1.Cells.Item(38, 6).Value = "🟡 Observer"
$ws1.Cells.Item(38, 7).Value = "👀 À surveiller"

# Row 39: TOTALENERGIES MARKETING SN (TTLS)
$ws1.Cells.Item(39, 1).Value = "TOTALENERGIES MARKETING SN (TTLS)"
$ws1.Cells.Item(39, 2).Value = 1
$ws1.Cells.Item(39, 3).Value = 1
$ws1.Cells.Item(39, 4).Value = -0.62
$ws1.Cells.Item(39, 5).Value = 3.02
$ws1.Cells.Item(39, 6).Value = "🟡 Observer"
$ws1.Cells.Item(39, 7).Value = "👀 À surveiller"

# Row 40: BANK OF AFRICA NG (BOAN)
$ws1.Cells.Item(40, 1).Value = "BANK OF AFRICA NG (BOAN)"
$ws1.Cells.Item(40, 2).Value = 1
$ws1.Cells.Item(40, 3).Value = 1
$ws1.Cells.Item(40, 4).Value = -1.24
$ws1.Cells.Item(40, 5).Value = 6
$ws1.Cells.Item(40, 6).Value = "🟡 Observer"
$ws1.Cells.Item(40, 7).Value = "👀 À surveiller"

# Row 41: VIVO ENERGY CI (SHEC)
$ws1.Cells.Item(41, 1).Value = "VIVO ENERGY CI (SHEC)"
$ws1.Cells.Item(41, 2).Value = 1
$ws1.Cells.Item(41, 3).Value = 2
$ws1.Cells.Item(41, 4).Value = -2.57
$ws1.Cells.Item(41, 5).Value = 4.43
$ws1.Cells.Item(41, 6).Value = "🟡 Observer"
$ws1.Cells.Item(41, 7).Value = "👀 À surveiller"

# Row 42: SMB CI (SMBC)
$ws1.Cells.Item(42, 1).Value = "SMB CI (SMBC)"
$ws1.Cells.Item(42, 2).Value = 1
$ws1.Cells.Item(42, 3).Value = 1
$ws1.Cells.Item(42, 4).Value = -2.93
$ws1.Cells.Item(42, 5).Value = 3.28
$ws1.Cells.Item(42, 6).Value = "🟡 Observer"
$ws1.Cells.Item(42, 7).Value = "👀 À surveiller"

# Row 43: BANK OF AFRICA ML (BOAM)
$ws1.Cells.Item(43, 1).Value = "BANK OF AFRICA ML (BOAM)"
$ws1.Cells.Item(43, 2).Value = 0
$ws1.Cells.Item(43, 3).Value = 1
$ws1.Cells.Item(43, 4).Value = -2.93
$ws1.Cells.Item(43, 5).Value = -2.93
$ws1.Cells.Item(43, 6).Value = "🟡 Observer"
$ws1.Cells.Item(43, 7).Value = "➖ Neutre"

# Row 44: SODE CI (SDCC)
$ws1.Cells.Item(44, 1).Value = "SODE CI (SDCC)"
$ws1.Cells.Item(44, 2).Value = 1
$ws1.Cells.Item(44, 3).Value = 1
$ws1.Cells.Item(44, 4).Value = -2.97
$ws1.Cells.Item(44, 5).Value = 2.43
$ws1.Cells.Item(44, 6).Value = "🟡 Observer"
$ws1.Cells.Item(44, 7).Value = "👀 À surveiller"

# Row 45: SONATEL SN (SNTS)
$ws1.Cells.Item(45, 1).Value = "SONATEL SN (SNTS)"
$ws1.Cells.Item(45, 2).Value = 0
$ws1.Cells.Item(45, 3).Value = 1
$ws1.Cells.Item(45, 4).Value = -4.23
$ws1.Cells.Item(45, 5).Value = -4.23
$ws1.Cells.Item(45, 6).Value = "🟡 Observer"
$ws1.Cells.Item(45, 7).Value = "➖ Neutre"

# Row 46: BANK OF AFRICA BF (BOABF)
$ws1.Cells.Item(46, 1).Value = "BANK OF AFRICA BF (BOABF)"
$ws1.Cells.Item(46, 2).Value = 0
$ws1.Cells.Item(46, 3).Value = 1
$ws1.Cells.Item(46, 4).Value = -5.33
$ws1.Cells.Item(46, 5).Value = -5.33
$ws1.Cells.Item(46, 6).Value = "🟡 Observer"
$ws1.Cells.Item(46, 7).Value = "➖ Neutre"

# Row 47: ECOBANK TRANS. INCORP. TG (ETIT)
$ws1.Cells.Item(47, 1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Cells.Item(47, 2).Value = 0
$ws1.Cells.Item(47, 3).Value = 1
$ws1.Cells.Item(47, 4).Value = -5.56
$ws1.Cells.Item(47, 5).Value = -5.56
$ws1.Cells.Item(47, 6).Value = "🟡 Observer"
$ws1.Cells.Item(47, 7).Value = "➖ Neutre"

# Row 48: SOLIBRA CI (SLBC)
$ws1.Cells.Item(48, 1).Value = "SOLIBRA CI (SLBC)"
$ws1.Cells.Item(48, 2).Value = 0
$ws1.Cells.Item(48, 3).Value = 1
$ws1.Cells.Item(48, 4).Value = -6.67
$ws1.Cells.Item(48, 5).Value = -6.67
$ws1.Cells.Item(48, 6).Value = "🟡 Observer"
$ws1.Cells.Item(48, 7).Value = "➖ Neutre"

# Row 49: AFRICA GLOBAL LOGISTICS CI (SDSC)
$ws1.Cells.Item(49, 1).Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws1.Cells.Item(49, 2).Value = 0
$ws1.Cells.Item(49, 3).Value = 2
$ws1.Cells.Item(49, 4).Value = -11.41
$ws1.Cells.Item(49, 5).Value = -7.42
$ws1.Cells.Item(49, 6).Value = "🟡 Observer"
$ws1.Cells.Item(49, 7).Value = "➖ Neutre"

# Row 50: SERVAIR ABIDJAN CI (ABJC)
$ws1.Cells.Item(50, 1).Value = "SERVAIR ABIDJAN CI (ABJC)"
$ws1.Cells.Item(50, 2).Value = 1
$ws1.Cells.Item(50, 3).Value = 3
$ws1.Cells.Item(50, 4).Value = -15.85
$ws1.Cells.Item(50, 5).Value = 4.24
$ws1.Cells.Item(50, 6).Value = "🔴 Vente"
$ws1.Cells.Item(50, 7).Value = "⚠️ Risque de décrochage"

# Row 51: FILTISAC CI (FTSC)
$ws1.Cells.Item(51, 1).Value = "FILTISAC CI (FTSC)"
$ws1.Cells.Item(51, 2).Value = 0
$ws1.Cells.Item(51, 3).Value = 4
$ws1.Cells.Item(51, 4).Value = -24.07
$ws1.Cells.Item(51, 5).Value = -7.4
$ws1.Cells.Item(51, 6).Value = "🔴 Vente"
$ws1.Cells.Item(51, 7).Value = "⚠️ Risque de décrochage"

# --- Sheet "Top_YTD": update rows 2-11 (cols A-B) ---
$ws2 = $wb.Worksheets.Item("Top_YTD")

# Row 2: BRVM - SERVICES PUBLICS
$ws2.Cells.Item(2, 1).Value = "BRVM - SERVICES PUBLICS"
$ws2.Cells.Item(2, 2).Value = 191352943.6

# Row 3: SAFCA CI
$ws2.Cells.Item(3, 1).Value = "SAFCA CI"
$ws2.Cells.Item(3, 2).Value = 3076463.36

# Row 4: CFAO MOTORS CI
$ws2.Cells.Item(4, 1).Value = "CFAO MOTORS CI"
$ws2.Cells.Item(4, 2).Value = 2795364.91

# Row 5: BRVM - AUTRES SECTEURS
$ws2.Cells.Item(5, 1).Value = "BRVM - AUTRES SECTEURS"
$ws2.Cells.Item(5, 2).Value = 2602501.54

# Row 6: UNIWAX CI
$ws2.Cells.Item(6, 1).Value = "UNIWAX CI"
$ws2.Cells.Item(6, 2).Value = 1609594.1

# Row 7: NEI-CEDA CI
$ws2.Cells.Item(7, 1).Value = "NEI-CEDA CI"
$ws2.Cells.Item(7, 2).Value = 1575099.47

# Row 8: SETAO CI
$ws2.Cells.Item(8, 1).Value = "SETAO CI"
$ws2.Cells.Item(8, 2).Value = 1105543.52

# Row 9: AIR LIQUIDE CI
$ws2.Cells.Item(9, 1).Value = "AIR LIQUIDE CI"
$ws2.Cells.Item(9, 2).Value = 1090484.58

# Row 10: BRVM - TRANSPORT
$ws2.Cells.Item(10, 1).Value = "BRVM - TRANSPORT"
$ws2.Cells.Item(10, 2).Value = 232634.54

# Row 11: BRVM - DISTRIBUTION
$ws2.Cells.Item(11, 1).Value = "BRVM - DISTRIBUTION"
$ws2.Cells.Item(11, 2).Value = 225256.48

